$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03786566666666667
$ws.Range("H2").Value = 0.113597
$ws.Range("I2").Value = 0.02729193434771035
$ws.Range("J2").Value = 0.02729193434771035
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3252056666666667
$ws.Range("N2").Value = 0.975617
$ws.Range("O2").Value = 0.0158278498560244
$ws.Range("P2").Value = 0.0158278498560244
$ws.Range("Q2").Value = 0.01231412937211111
$ws.Range("R2").Value = 0.110827164349
$ws.Range("S2").Value = 0.0004319726391360346
$ws.Range("T2").Value = 0.0004319726391360346

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03786566666666667
$ws.Range("H3").Value = 0.113597
$ws.Range("I3").Value = 0.02729193434771035
$ws.Range("J3").Value = 0.02729193434771035
$ws.Range("O3").Value = 0.8133441666880411
$ws.Range("P3").Value = 0.8133441666880411
$ws.Range("Q3").Value = 0.6327849571327778
$ws.Range("R3").Value = 5.695064614195
$ws.Range("S3").Value = 0.0221977355993432
$ws.Range("T3").Value = 0.0221977355993432

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03786566666666667
$ws.Range("H4").Value = 0.113597
$ws.Range("I4").Value = 0.02729193434771035
$ws.Range("J4").Value = 0.02729193434771035
$ws.Range("M4").Value = 3.509903666666667
$ws.Range("N4").Value = 10.529711
$ws.Range("O4").Value = 0.1708279834559346
$ws.Range("P4").Value = 0.1708279834559346
$ws.Range("Q4").Value = 0.1329048422741111
$ws.Range("R4").Value = 1.196143580467
$ws.Range("S4").Value = 0.004662226109231117
$ws.Range("T4").Value = 0.004662226109231117

# Row 5
$ws.Range("I5").Value = 0.7739041374319726
$ws.Range("J5").Value = 0.7739041374319726
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3252056666666667
$ws.Range("N5").Value = 0.975617
$ws.Range("O5").Value = 0.0158278498560244
$ws.Range("P5").Value = 0.0158278498560244
$ws.Range("Q5").Value = 0.3491857905172222
$ws.Range("R5").Value = 3.142672114655
$ws.Range("S5").Value = 0.01224923849022933
$ws.Range("T5").Value = 0.01224923849022933

# Row 6
$ws.Range("I6").Value = 0.7739041374319726
$ws.Range("J6").Value = 0.7739041374319726
$ws.Range("O6").Value = 0.8133441666880411
$ws.Range("P6").Value = 0.8133441666880411
$ws.Range("S6").Value = 0.629450415756035
$ws.Range("T6").Value = 0.629450415756035

# Row 7
$ws.Range("I7").Value = 0.7739041374319726
$ws.Range("J7").Value = 0.7739041374319726
$ws.Range("M7").Value = 3.509903666666667
$ws.Range("N7").Value = 10.529711
$ws.Range("O7").Value = 0.1708279834559346
$ws.Range("P7").Value = 0.1708279834559346
$ws.Range("Q7").Value = 3.768718113207222
$ws.Range("R7").Value = 33.918463018865
$ws.Range("S7").Value = 0.1322044831857083
$ws.Range("T7").Value = 0.1322044831857083

# Row 8
$ws.Range("G8").Value = 0.2758266666666667
$ws.Range("H8").Value = 0.82748
$ws.Range("I8").Value = 0.1988039282203171
$ws.Range("J8").Value = 0.1988039282203171
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3252056666666667
$ws.Range("N8").Value = 0.975617
$ws.Range("O8").Value = 0.0158278498560244
$ws.Range("P8").Value = 0.0158278498560244
$ws.Range("Q8").Value = 0.08970039501777778
$ws.Range("R8").Value = 0.80730355516
$ws.Range("S8").Value = 0.003146638726659031
$ws.Range("T8").Value = 0.003146638726659031

# Row 9
$ws.Range("G9").Value = 0.2758266666666667
$ws.Range("H9").Value = 0.82748
$ws.Range("I9").Value = 0.1988039282203171
$ws.Range("J9").Value = 0.1988039282203171
$ws.Range("O9").Value = 0.8133441666880411
$ws.Range("P9").Value = 0.8133441666880411
$ws.Range("Q9").Value = 4.609425392644445
$ws.Range("R9").Value = 41.4848285338
$ws.Range("S9").Value = 0.161696015332663
$ws.Range("T9").Value = 0.161696015332663

# Row 10
$ws.Range("G10").Value = 0.2758266666666667
$ws.Range("H10").Value = 0.82748
$ws.Range("I10").Value = 0.1988039282203171
$ws.Range("J10").Value = 0.1988039282203171
$ws.Range("M10").Value = 3.509903666666667
$ws.Range("N10").Value = 10.529711
$ws.Range("O10").Value = 0.1708279834559346
$ws.Range("P10").Value = 0.1708279834559346
$ws.Range("Q10").Value = 0.9681250286977778
$ws.Range("R10").Value = 8.71312525828
$ws.Range("S10").Value = 0.03396127416099513
$ws.Range("T10").Value = 0.03396127416099513

Write-Host "Update complete"